$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. "1.0000",
# "0.06500") keep their exact literal formatting instead of being coerced
# into numbers by Excel. Style is reset to Normal afterwards so no stray
# cell formatting is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '30.424.09'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '1.871.02'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '246.37'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.4743'
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('D8').Value = '0.2922'
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('D9').Value = '0.06500'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').Value = '22.28'
$ws.Range('E10').Value = '  +7.59%  '
$ws.Range('D11').Value = '0.07716'
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').Value = '97.77'
$ws.Range('E12').Value = '  +3.42%  '
$ws.Range('D13').Value = '0.7386'
$ws.Range('E13').Value = '  +6.18%  '
$ws.Range('D14').Value = '1.872.60'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '5.140'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('D16').Value = '273.31'
$ws.Range('E16').Value = '  +2.13%  '
$ws.Range('D17').Value = '30.408.26'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').Value = '13.40'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '0.000007542'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = '2.116.71'
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '5.219'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').Value = '6.172'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').Value = '9.293'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').Value = '163.48'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '18.83'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').Value = '1.925'
$ws.Range('E28').Value = '  +1.60%  '
$ws.Range('D29').Value = '0.1005'
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('D30').Value = '1.367'
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('D31').Value = '1.505'
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('D32').Value = '4.288'
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('D33').Value = '4.135'
$ws.Range('E33').Value = '  +3.49%  '
$ws.Range('D34').Value = '0.04824'
$ws.Range('E34').Value = '  +2.98%  '
$ws.Range('E35').Value = '  +1.03%  '
$ws.Range('D36').Value = '0.6942'
$ws.Range('E36').Value = '  +1.12%  '
$ws.Range('D37').Value = '0.9996'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '2.717'
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('D39').Value = '0.01853'
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('D40').Value = '2.743'
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = '6.302'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = '73.09'
$ws.Range('E42').Value = '  +4.08%  '
$ws.Range('D43').Value = '1.965'
$ws.Range('E43').Value = '  +4.24%  '
$ws.Range('D44').Value = '0.4194'
$ws.Range('E44').Value = '  +3.46%  '
$ws.Range('D45').Value = '0.9999'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '0.8333'
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('D47').Value = '102.08'
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').Value = '9.231'
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Value = '7.011'
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '35.46'
$ws.Range('E50').Value = '  +2.98%  '
$ws.Range('D51').Value = '924.10'
$ws.Range('E51').Value = '  -0.60%  '

$dRange.Style = "Normal"

